# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are numeric-looking text (e.g. "7.20", "26.705.84")
# Force text storage so trailing zeros / thousand-dot formatting survive,
# matching the original inline-string cells (mirrors formatting the column as Text).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "26.705.84"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "1.531.56"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "205.86"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.484"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "21.35"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.0579"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "1.747.33"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "1.542.13"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.506"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "61.34"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "26.675.32"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "212.57"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "7.20"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "4.01"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "152.29"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "6.54"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "14.85"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.10"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "0.0452"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "1.355.64"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.946"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.522"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.797"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "5.67"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.993"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "62.46"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "1.663.19"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "85.48"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.0506"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.0943"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.0₇0955"

# Coin name / link / 1h-volume columns are plain text already - no coercion needed
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  -2.31%  "
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("E23").Value = "  -3.01%  "
$ws.Range("E24").Value = "  -3.84%  "
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  -3.44%  "
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("E31").Value = "  -2.13%  "
$ws.Range("E32").Value = "  +2.80%  "
$ws.Range("E33").Value = "  -2.03%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  -3.67%  "
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E41").Value = "  -1.74%  "
$ws.Range("E42").Value = "  +5.26%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  +2.78%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("E51").Value = "  -1.71%  "
